# Events.xlsx - "Finialized transition to new day"
#
# The e006a/e006b "Retrofit Period" entries are being split into three
# distinct events (e006a, e006b, e006c):
#   - Row 8  (e006a): trim the stray trailing blank line from the body text
#             and shrink the row to fit (150 -> 135).
#   - Row 9  (e006b): left as-is (Crew Training body text unchanged).
#   - Row 10 (was a duplicate "e006b"): relabeled to the new "e006c" event,
#             with its body text's bold heading updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: e006a Retrofit Period -----------------------------------------
# Remove the single trailing newline left at the end of the cell text.
$b8 = $ws.Range("B8").Value2
$ws.Range("B8").Value2 = $b8.Substring(0, $b8.Length - 1)
$ws.Rows.Item(8).RowHeight = 135

# --- Row 10: new e006c Retrofit Period - Gyrostabilizer -------------------
$ws.Range("A10").Value2 = "e006c"
$b10 = $ws.Range("B10").Value2
$ws.Range("B10").Value2 = $b10.Replace("e006b Retrofit Period - Gyrostabilizer", "e006c Retrofit Period - Gyrostabilizer")

# --- Selection/view housekeeping -------------------------------------------
$ws.Range("B9").Select() | Out-Null
